$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last data row (row 31) onto the new row 32
# so the new row matches the shared fill/alignment styles used by the
# rest of the comparison table.
$ws.Range("A31:I31").Copy()
$ws.Range("A32:I32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add a new data row documenting the "userLoginSource" / "string" field
# pair, with the usual A=D and B=E check formulas in F:G. Columns C, H
# and I stay empty, like the rest of the table.
$ws.Range("A32").Value = "userLoginSource"
$ws.Range("B32").Value = "string"
$ws.Range("D32").Value = "userLoginSource"
$ws.Range("E32").Value = "string"
$ws.Range("F32").Formula = "=A32=D32"
$ws.Range("G32").Formula = "=B32=E32"
$ws.Range("C32").ClearContents()
$ws.Range("H32").ClearContents()
$ws.Range("I32").ClearContents()

# Move the active selection to C27, as recorded in the saved workbook.
$ws.Range("C27").Select()
